$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 946:948, shifting existing rows 946-994 down to 949-997
$ws.Range('A946:R948').Insert()

# Fill in the 3 newly inserted rows with the new weekly data
$ws.Cells.Item(946, 1).Value = 11
$ws.Cells.Item(946, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(946, 3).Value = 'Bíobío'
$ws.Cells.Item(946, 4).Value = 45267
$ws.Cells.Item(946, 5).Value = 8
$ws.Cells.Item(946, 6).Value = 100112004
$ws.Cells.Item(946, 7).Value = 'Cebolla'
$ws.Cells.Item(946, 8).Value = 'Morada(o)'
$ws.Cells.Item(946, 9).Value = 'Primera'
$ws.Cells.Item(946, 10).Value = 200
$ws.Cells.Item(946, 11).Value = 10000
$ws.Cells.Item(946, 12).Value = 11000
$ws.Cells.Item(946, 13).Value = 10500
$ws.Cells.Item(946, 14).Value = '$/malla 18 kilos'
$ws.Cells.Item(946, 15).Value = 'Perú'
$ws.Cells.Item(946, 16).Value = 583
$ws.Cells.Item(946, 17).Value = 18
$ws.Cells.Item(946, 18).Value = 'Hortaliza'

$ws.Cells.Item(947, 1).Value = 11
$ws.Cells.Item(947, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(947, 3).Value = 'Bíobío'
$ws.Cells.Item(947, 4).Value = 45267
$ws.Cells.Item(947, 5).Value = 8
$ws.Cells.Item(947, 6).Value = 100112004
$ws.Cells.Item(947, 7).Value = 'Cebolla'
$ws.Cells.Item(947, 8).Value = 'Sin especificar'
$ws.Cells.Item(947, 9).Value = '1a (cosecha)'
$ws.Cells.Item(947, 10).Value = 400
$ws.Cells.Item(947, 11).Value = 11000
$ws.Cells.Item(947, 12).Value = 12000
$ws.Cells.Item(947, 13).Value = 11500
$ws.Cells.Item(947, 14).Value = '$/malla 16 kilos'
$ws.Cells.Item(947, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(947, 16).Value = 719
$ws.Cells.Item(947, 17).Value = 16
$ws.Cells.Item(947, 18).Value = 'Hortaliza'

$ws.Cells.Item(948, 1).Value = 11
$ws.Cells.Item(948, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(948, 3).Value = 'Bíobío'
$ws.Cells.Item(948, 4).Value = 45267
$ws.Cells.Item(948, 5).Value = 8
$ws.Cells.Item(948, 6).Value = 100112004
$ws.Cells.Item(948, 7).Value = 'Cebolla'
$ws.Cells.Item(948, 8).Value = 'Sin especificar'
$ws.Cells.Item(948, 9).Value = '2a (cosecha)'
$ws.Cells.Item(948, 10).Value = 200
$ws.Cells.Item(948, 11).Value = 9000
$ws.Cells.Item(948, 12).Value = 9000
$ws.Cells.Item(948, 13).Value = 9000
$ws.Cells.Item(948, 14).Value = '$/malla 16 kilos'
$ws.Cells.Item(948, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(948, 16).Value = 562
$ws.Cells.Item(948, 17).Value = 16
$ws.Cells.Item(948, 18).Value = 'Hortaliza'
